# feat: Add total price calculation and styling
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Sales Report"

# New "total price" column header + values (price * quantity)
$ws.Range("D1").Value = "ราคารวม (บาท)"
$ws.Range("D2").Value = 3000
$ws.Range("D3").Value = 4750
$ws.Range("D4").Value = 1200

# Column widths: A (item name) wider, D (total) a bit narrower
# ColumnWidth is in characters; Excel stores width in the sheet as
# characters + ~5/6 padding, so bias the input so the stored width
# lands exactly on 20 / 15.
$ws.Columns.Item(1).ColumnWidth = 19.16666667
$ws.Columns.Item(4).ColumnWidth = 14.16666667

# Header row styling: bold white text on a solid blue fill
$headerRange = $ws.Range("A1:D1")
$headerRange.Interior.Color = 12419407
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215

Write-Host "Done"
